# Update the workers' account-statement table (rows 16-22) with the new
# roster of debtors/amounts, replacing the previous data set, and bump the
# overdue-value for OSCAR ADAN MARRUGO MARTINEZ.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Doc type (B), Doc number (C), Name (D), Period (E), Salario (F), Valor Mora (G)
$rows = @(
    @{ Row = 16; DocType = "CC"; DocNum = "8852039";    Name = "JHONNI FRANCISCO PADILLA PACHECO"; Period = "2507"; Salario = 56940; Valor = 1423500 },
    @{ Row = 17; DocType = "CC"; DocNum = "9159074";     Name = "ADOLFREDO ISIDRO TERAN ZUÑIGA";     Period = "2507"; Salario = 56940; Valor = 1423500 },
    @{ Row = 18; DocType = "CC"; DocNum = "1026561058";  Name = "JUAN DAMASO ZABALETA FLOREZ";       Period = "2507"; Salario = 56940; Valor = 1423500 },
    @{ Row = 19; DocType = "CC"; DocNum = "73432563";    Name = "EDUARDO ENRIQUE SALCEDO CARO";      Period = "2507"; Salario = 56940; Valor = 1423500 },
    @{ Row = 20; DocType = "CC"; DocNum = "1143361541";  Name = "YESSICA PATRICIA CARMONA HERRERA";  Period = "2507"; Salario = 56940; Valor = 1423500 },
    @{ Row = 21; DocType = "CC"; DocNum = "1143394318";  Name = "OSCAR ADAN MARRUGO MARTINEZ";       Period = "2507"; Salario = 56940; Valor = 4200000 },
    @{ Row = 22; DocType = "CC"; DocNum = "73008467";    Name = "AURELIO MANUEL YEPES GUTIERREZ";    Period = "2507"; Salario = 56940; Valor = 1423500 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 2).Value = $r.DocType
    $ws.Cells.Item($row, 3).Value = $r.DocNum
    $ws.Cells.Item($row, 4).Value = $r.Name
    $ws.Cells.Item($row, 5).Value = $r.Period
    $ws.Cells.Item($row, 6).Value = $r.Salario
    $ws.Cells.Item($row, 7).Value = $r.Valor
}
